$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Excel auto-converts a typed "TRUE"/"FALSE" into the Boolean literal, even in
# a text-formatted cell, unless it is entered with a leading apostrophe (quote
# prefix) to force literal text. Stage the literal text value in a scratch
# cell using the apostrophe trick, then copy/paste-special (values only) it
# onto each target cell so the destination keeps its own existing style
# (text number format) instead of picking up the scratch cell's "quote
# prefix" style.
$helper = $ws.Cells.Item(50, 50)
$helper.Value = "'TRUE"
$helper.Copy()

# Replace the TRUE() boolean formulas in column E (rows 2-13) with the literal
# text string "TRUE" so the cells store a shared string instead of a numeric
# formula result.
for ($row = 2; $row -le 13; $row++) {
    $cell = $ws.Cells.Item($row, 5)
    $cell.PasteSpecial(-4163)  # xlPasteValues
}

$helper.ClearContents()
$excel.CutCopyMode = $false

# Move the active selection from the whole E column to cell E2, matching the
# updated sheet view state.
$ws.Range("E2").Select()
